$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in G1 / H1, matching the style used by the existing headers (F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Fill the new "Elapsed Time" / "CPU" columns for every data row
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 7).Value = 1.669922641383406
    $ws.Cells.Item($r, 8).Value = 0.97
}
